$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 237 (B and D values changed)
$ws.Range("B237").Value = 177772000000000
$ws.Range("D237").Value = 136015302218.8217

# Copy formatting of row 237 (A column date style) down to new rows 238:240
$ws.Range("A237").Copy()
$ws.Range("A238:A240").PasteSpecial(-4122)

# Add new row 238
$ws.Range("A238").Value = 45108
$ws.Range("B238").Value = 176788000000000
$ws.Range("C238").Value = 0.0007645259938837921
$ws.Range("D238").Value = 135159021406.7278

# Add new row 239
$ws.Range("A239").Value = 45139
$ws.Range("B239").Value = 174322000000000
$ws.Range("C239").Value = 0.0007641637755803825
$ws.Range("D239").Value = 133210557686.7234

# Add new row 240
$ws.Range("A240").Value = 45170
$ws.Range("B240").Value = 173950000000000
$ws.Range("C240").Value = 0.0007641637755803825
$ws.Range("D240").Value = 132926288762.2075
